$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: test objective updated ---
$ws.Range("C2").Value = "Verify successful company creation with all mandatory fields."

# --- Row 3: test objective updated (was placeholder "Verify Tc") ---
$ws.Range("C3").Value = "Verify that company name field can be edited and updated value is displayed correctly."

# --- New rows 17-24 (TC_016 .. TC_023) ---
# Seed each new row's formatting (styles + row height) from row 16, the last existing row.
$ws.Range("A16:I16").Copy()
$ws.Range("A17:I24").PasteSpecial(-4122)
for ($r = 17; $r -le 24; $r++) {
    $ws.Rows($r).RowHeight = 80
}
# Serial numbers in column A are plain text (like "01".."15" above them), so force
# text storage before writing "16".."23" - otherwise they'd be auto-parsed as numbers.
$ws.Range("A17:A24").NumberFormat = "@"

# Row 20's "Test Status" cell (I20) becomes the new FAILED row, inheriting the
# red/pink FAILED formatting that I4 currently has. Copy that formatting over
# BEFORE row 4 itself is reset back to "Not Run" below.
$ws.Range("I4").Copy()
$ws.Range("I20").PasteSpecial(-4122)

# --- Row 4: result flipped from FAILED back to Not Run ---
# Reuse the existing "Not Run" formatting (style 3) from I2 so no new style is minted.
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("H4").Value = "Test not executed"
$ws.Range("I4").Value = "Not Run"

$boilerPre = "User is on the company field editing page"
$boilerData = "Company Field Editing test data"
$boilerSteps = "1. Navigate to company field editing page`n2. Perform required actions`n3. Verify expected behavior"
$boilerExpected = "Company Field Editing functionality should work as expected"
$notExecuted = "Test not executed"
$notRun = "Not Run"
$failedText = "Test failed - actual behavior did not match expected result"
$failedStatus = "FAILED"

$rows = @(
    @{ Row = 17; Serial = "16"; Tc = "TC_016"; Objective = "Verify comprehensive company fields editing using direct locators approach."; Failed = $false },
    @{ Row = 18; Serial = "17"; Tc = "TC_017"; Objective = "Verify individual company name field editing with current value."; Failed = $false },
    @{ Row = 19; Serial = "18"; Tc = "TC_018"; Objective = "Verify individual website field editing with current value."; Failed = $false },
    @{ Row = 20; Serial = "19"; Tc = "TC_019"; Objective = "Verify individual industry field editing with current value."; Failed = $true },
    @{ Row = 21; Serial = "20"; Tc = "TC_020"; Objective = "Verify individual HQ in JPN field editing with current value."; Failed = $false },
    @{ Row = 22; Serial = "21"; Tc = "TC_021"; Objective = "Verify individual Global HQ field editing with current value."; Failed = $false },
    @{ Row = 23; Serial = "22"; Tc = "TC_022"; Objective = "Verify individual Country of origin field editing with current value."; Failed = $false },
    @{ Row = 24; Serial = "23"; Tc = "TC_023"; Objective = "Verify individual Company address field editing with current value."; Failed = $false }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Serial
    $ws.Range("B$r").Value = $item.Tc
    $ws.Range("C$r").Value = $item.Objective
    $ws.Range("D$r").Value = $boilerPre
    $ws.Range("E$r").Value = $boilerData
    $ws.Range("F$r").Value = $boilerSteps
    $ws.Range("G$r").Value = $boilerExpected

    if ($item.Failed) {
        $ws.Range("H$r").Value = $failedText
        $ws.Range("I$r").Value = $failedStatus
    } else {
        $ws.Range("H$r").Value = $notExecuted
        $ws.Range("I$r").Value = $notRun
    }
}
